$d = $word.ActiveDocument

# Locate the standalone "Alternate" run (whole-word, case-sensitive) that
# needs to become "Alternative". Other occurrences of "Alternative"
# elsewhere in the document are left untouched because they do not match
# "Alternate" as a whole word.
$find = $d.Content
[void]$find.Find.Execute("Alternate", $true, $true, $false, $false, $false, `
                          $true, 1, $false, "", 0)

if ($find.Find.Found) {
    $wordStart = $find.Start
    $wordEnd = $find.End

    # Keep "Alternat" (first 8 chars) as-is; replace the trailing "e"
    # with "ive" in a *separate* run so the final text reads
    # "Alternat" + "ive" = "Alternative", matching the two-run split in
    # the target markup (both runs keep identical Bold/BCs/Underline
    # formatting).
    $tailStart = $wordStart + 8
    $tail = $d.Range($tailStart, $wordEnd)

    # Toggling the underline off and back on forces the engine to keep
    # this as its own run instead of silently re-merging it with the
    # preceding "Alternat" run once the formatting matches again.
    $tail.Font.Underline = 0
    $tail.Text = "ive"

    $newRun = $d.Range($tailStart, $tailStart + 3)
    $newRun.Font.Bold = $true
    $newRun.Font.Underline = 1
}
